$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ============ 1. Resize table A1:D6 -> A1:F7 ============
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F7"))

# ============ 2. Column widths ============
$ws.Columns.Item(4).ColumnWidth = 13.5
$ws.Columns.Item(5).ColumnWidth = 63.65
$ws.Columns.Item(6).ColumnWidth = 30.2

# ============ 3. Header row (row 1) ============
$ws.Range("C1").Value = "Goals"
$ws.Range("D1").Value = "Used methods"
$ws.Range("E1").Value = "Results"
$ws.Range("F1").Value = "Why important (one sentence)"

# ============ 4. Row 2 data (C2,D2,E2) ============
$ws.Range("C2").Value = 'Remove the good ideas behind MDE, its limitation and the bad part of MDE'
$ws.Range("D2").Value = 'Analysis of MDE and its limitations through analysis'
$ws.Range("E2").Value = 'A desire to reduce the time to market for software systems, their continuously growing complexity, and ongoing demand for higher software quality was an incentive to look for some automation tools. As a result, the concept of MDE appeared; A straightforward but not the smartest approach to implement some system is to have some log with requirements and present architecture as some figures with available interfaces. Based on it, developers tend directly writing source code. The list of bugs is kept in some issue-tracking products. Taking into account that by such an approach the requirements are not mapped to architecture and source code, then it becomes extremely difficult to reflect this match during the implementation stage. Gradually it might become even worse as soon as the code base begins to age and produce even more bugs. The solution to overcome those problems seems obvious. If one has an overarching single data source that would link to each other all the stages of software development, it would help keep all aspects of development interconnected to each other;'

# ============ 5. Row 6 data (C6,D6,E6,F6) ============
$ws.Range("C6").Value = 'Give an overview of what the MDE is, its origination roots, and its current state; specifically how it can be used in modern software systems'
$ws.Range("D6").Value = 'Course defitions of a model, mde, model transformation'
$ws.Range("E6").Value = 'Model is a graph-oriented structure to present a certain domain in its simplified but coherent representation which adheres to another graph structure named metamodel; MDE is an approach that relies on three components that are closely related to each other: automation, DSL, and commonly-agreed standards; MDE is as a technique to create software products by considering models and their elements as the most important construction blocks; '
$ws.Range("F6").Value = 'Definitions of model, mde and its key componetns'

# ============ 6. Row 7 (new row, full citation) ============
$ws.Range("A7").Value = '"@article{bezivin2004search,
  title={In search of a basic principle for model driven engineering},
  author={B{\''e}zivin, Jean},
  journal={Novatica Journal, Special Issue},
  volume={5},
  number={2},
  pages={21--24},
  year={2004},
  publisher={Citeseer}
}"'
$ws.Range("B7").Value = 'In search of a basic principle for model driven engineering'
$ws.Range("C7").Value = 'find fundamental ideas of MDE'
$ws.Range("D7").Value = 'Similar to the term "everything is an object" define a statement "everything is a model" in the scope of MDA to look for the essential trates of MDE'
$ws.Range("E7").Value = 'MDA is regarded as a practical implementation of MDE that was created with the use of OMG standards;
a Switch of the paradigm to viewing not an object as a central element of the systems but model and its elements. This allows us to consider models not just as a form of representing documentation but as a complete space to drive software product lines; 
After the introduction of MDA by OMG Group, the transition from code-based to model-based software happened. This, in turn, induced the appearance of many languages used to specify a certain domain. The variety of different meta-languages urged the need to create a unified framework that all meta-languages could conform to and, thus, to make them interchangeable. That is how emerged MOF - a unified modeling language for all meta-models;'
$ws.Range("F7").Value = 'What is MDA'

# ============ 7. Formatting: copy standard row styles onto new cells ============
# E/F columns for rows 3,4,5,7 use the standard "last-col" style (border6, center/center/wrap)
$ws.Range("D4").Copy() | Out-Null
$ws.Range("E3:F5").PasteSpecial(-4122) | Out-Null
$ws.Range("E7:F7").PasteSpecial(-4122) | Out-Null
# Row 7 A:D need the same per-column styles as the other data rows (copy from row 4)
$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A7:D7").PasteSpecial(-4122) | Out-Null

# ============ 8. New style for row2 E/F: left+bottom thin border, center/center/wrap ============
$scratchA = $ws.Range("Z100")
$scratchA.WrapText = $true
$scratchA.HorizontalAlignment = -4108
$scratchA.VerticalAlignment = -4108
$scratchA.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$scratchA.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$scratchA.Copy() | Out-Null
$ws.Range("E2:F2").PasteSpecial(-4122) | Out-Null
$scratchA.Clear() | Out-Null

# ============ 9. New style for row6 E/F: left+top thin border, center/center/wrap ============
$scratchB = $ws.Range("Z101")
$scratchB.WrapText = $true
$scratchB.HorizontalAlignment = -4108
$scratchB.VerticalAlignment = -4108
$scratchB.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$scratchB.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$scratchB.Copy() | Out-Null
$ws.Range("E6:F6").PasteSpecial(-4122) | Out-Null
$scratchB.Clear() | Out-Null

# ============ 10. Row heights ============
$ws.Rows.Item(1).RowHeight = 28.8
$ws.Rows.Item(7).RowHeight = 409.6

# ============ 11. Selection / view ============
$ws.Range("C3").Select() | Out-Null

Write-Host "done"
